$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching the style of the existing header row (column E)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamp values for rows 2-19 (unstyled, plain text like the rest of the data rows)
$timestamps = @(
    "2021-10-05 13:41:34.596026",
    "2021-10-05 13:41:34.596038",
    "2021-10-05 13:41:34.596042",
    "2021-10-05 13:41:34.596045",
    "2021-10-05 13:41:34.596048",
    "2021-10-05 13:41:34.596051",
    "2021-10-05 13:41:34.596054",
    "2021-10-05 13:41:34.596057",
    "2021-10-05 13:41:34.596060",
    "2021-10-05 13:41:34.596063",
    "2021-10-05 13:41:34.596066",
    "2021-10-05 13:41:34.596069",
    "2021-10-05 13:41:34.596072",
    "2021-10-05 13:41:34.596084",
    "2021-10-05 13:41:34.596087",
    "2021-10-05 13:41:34.596090",
    "2021-10-05 13:41:34.596093",
    "2021-10-05 13:41:34.596097"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
